$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits at the
#    very start of the document (before the "Eflows" heading).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Bold the whole "Take class wide DRH ... falls within range"
#    bullet (paragraph mark included, so both runs + the pPr/rPr
#    pick up <w:b/>).
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Take class wide DRH*falls within range*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Bold = 1
}

# ------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark in the middle of the "Decide as
#    a group ..." bullet, splitting the run right after "...whether
#    to in" / before "clude everything ...".
# ------------------------------------------------------------------
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Decide as a group on a threshold*") {
        $target2 = $p
        break
    }
}
if ($target2 -ne $null) {
    $r = $target2.Range
    $paraText = $r.Text
    $splitOffset = $paraText.IndexOf("clude everything")
    if ($splitOffset -ge 0) {
        $splitPos = $r.Start + $splitOffset
        $bmRange = $d.Range($splitPos, $splitPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)
    }
}
